# Restored from revision #7f2ef7bfd90e8c783b3e1682a036b9530583f613.TEST
# Author: admin. Type: SAVE.
#
# Change: cell E8 on the "Rules" sheet held the shared string "Good Morning";
# it is updated to the (new) string "Good Mornin".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "Good Mornin"
